$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "Retry" sheet selection (it is no longer the active/selected
#    tab once the new sheet is added - that happens automatically - but its
#    stored selection also changes from L5 to the whole used range A1:L8).
# ---------------------------------------------------------------------------
$wsRetry = $wb.Worksheets.Item("Retry")
$wsRetry.Range("A1:L8").Select()

# ---------------------------------------------------------------------------
# 2. Add the new "New" worksheet after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add($null, $lastSheet)
$wsNew.Name = "New"

# ---------------------------------------------------------------------------
# 3. Populate the header row (row 1).
#    Columns: B..P -> Basic, Spinning, Fire, Glowing, x2, x3, x4, x5, x6, x7,
#             x8, x9, x10, x11, x12
# ---------------------------------------------------------------------------
$headers = @("Basic","Spinning","Fire","Glowing","x2","x3","x4","x5","x6","x7","x8","x9","x10","x11","x12")
$col = 2
foreach ($h in $headers) {
    $cell = $wsNew.Cells.Item(1, $col)
    $cell.Value = $h
    $cell.Font.Bold = $true
    $col = $col + 1
}

# ---------------------------------------------------------------------------
# 4. Populate column A (row labels) for rows 2-8.
#    Note: matches the existing Fail/Retry sheets' row order (Blue comes
#    before Yellow).
# ---------------------------------------------------------------------------
$labels = @("Red","Green","Blue","Yellow","Purple","Orange","White")
$row = 2
foreach ($lbl in $labels) {
    $cell = $wsNew.Cells.Item($row, 1)
    $cell.Value = $lbl
    $cell.Font.Bold = $true
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 5. Populate the data grid (rows 2-8, columns B-O) with "X" in red font.
#    Some cells are intentionally left blank (still red-formatted) and some
#    are omitted entirely, matching the source data.
# ---------------------------------------------------------------------------
# Columns, 1-indexed: A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8 I=9 J=10 K=11 L=12 M=13 N=14 O=15 P=16
$filledCols = @{
    2  = @(2,3,4,5,6,7,8,9,10,11,12,13,14)        # row 2: B..N
    3  = @(2,3,4,5,6,7,8,9,10,11,13,14)            # row 3: B..K, M, N  (L blank)
    4  = @(2,3,4,5,6,7,8,9,10,11,12,13)            # row 4: B..M
    5  = @(2,3,4,5,6,7,8,9,10,11)                  # row 5: B..K (L blank)
    6  = @(2,3,4,5,6,7,8,9,10,11,12,14)            # row 6: B..L, N (M missing)
    7  = @(2,3,4,5,6,7,8,9,10,11,12,13,14)         # row 7: B..N
    8  = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15)      # row 8: B..O
}
$blankCols = @{
    3 = @(12)   # L3
    5 = @(12)   # L5
}

foreach ($r in 2..8) {
    foreach ($c in $filledCols[$r]) {
        $cell = $wsNew.Cells.Item($r, $c)
        $cell.Value = "X"
        $cell.Font.Color = 255
    }
    if ($blankCols.ContainsKey($r)) {
        foreach ($c in $blankCols[$r]) {
            $cell = $wsNew.Cells.Item($r, $c)
            $cell.Font.Color = 255
        }
    }
}

# ---------------------------------------------------------------------------
# 6. Set the selection/active cell on the new sheet and make sure it is the
#    active (tab-selected) sheet.
# ---------------------------------------------------------------------------
$wsNew.Range("M9").Select()
